$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching style of existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new columns
$values = @{
    2 = 8
    3 = 9
    4 = 8
    5 = 8
    6 = 7
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 9).Value = $val   # column I
    $ws.Cells.Item($row, 10).Value = $val  # column J
}
